$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) values
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 380
$wsExhibition.Range("F4").Value = 1589
$wsExhibition.Range("F5").Value = 15
$wsExhibition.Range("F6").Value = 23
$wsExhibition.Range("F9").Value = 63
$wsExhibition.Range("F10").Value = 468

# Sheet "全部类型" (sheet4): update column F (想去人数) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 380
$wsAll.Range("F3").Value = 0
$wsAll.Range("F6").Value = 0
$wsAll.Range("F7").Value = 408
$wsAll.Range("F8").Value = 141
$wsAll.Range("F9").Value = 63
$wsAll.Range("F10").Value = 0
